$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Edit 3 (done first, while the old "_GoBack" bookmark still exists):
#   "Oberflächen-Prototyp erstellen (ohne " + bookmark(_GoBack) + "größere Funktionalität)"
#   -> single run "Oberflächen-Prototyp erstellen (ohne größere Funktionalität)"
#   (the bookmark is removed entirely)
# -----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks("_GoBack")
    $oldPara = $oldBm.Range.Paragraphs(1)
    $oldParaStart = $oldPara.Range.Start

    # Drop the bookmark itself first.
    $oldBm.Delete()

    # Re-acquire the (now bookmark-free) paragraph via its original start
    # position, excluding the trailing paragraph mark, and force it into a
    # single run by writing a throwaway placeholder before writing the real
    # text back (a plain re-assignment of the identical text is a no-op and
    # would leave the existing run split untouched).
    $mergePara = $d.Range($oldParaStart, $oldParaStart).Paragraphs(1)
    $mergeRng = $d.Range($mergePara.Range.Start, $mergePara.Range.End - 1)
    $mergeRng.Text = "__PLACEHOLDER__"

    $mergePara2 = $d.Range($oldParaStart, $oldParaStart).Paragraphs(1)
    $mergeRng2 = $d.Range($mergePara2.Range.Start, $mergePara2.Range.End - 1)
    $mergeRng2.Text = "Oberflächen-Prototyp erstellen (ohne größere Funktionalität)"
}

# -----------------------------------------------------------------------
# Edit 2: "Oberfläche programmieren" -> two runs: "Oberfläche " / "erstellen"
# (no bookmark remains here; a temporary bookmark is used purely to force
#  the run split, then removed again)
# -----------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Oberfläche programmieren", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Oberfläche programmieren'"
}
$start2 = $rng2.Start
$splitPos2 = $start2 + 11  # length of "Oberfläche " (10 letters + trailing space)
$tailRng2 = $d.Range($splitPos2, $rng2.End)
if ($tailRng2.Text -ne "programmieren") {
    throw "Unexpected tail text at edit 2: '$($tailRng2.Text)'"
}
$tailRng2.Text = "erstellen"

$tmpMarkName = "ZZTempSplit"
$tmpBmRng2 = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add($tmpMarkName, $tmpBmRng2)
$d.Bookmarks($tmpMarkName).Delete()

# -----------------------------------------------------------------------
# Edit 1: "Arbeitspakete:" -> "Arbeitspaket" + bookmark(_GoBack) + ":"
# -----------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Arbeitspakete:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Arbeitspakete:'"
}
$start1 = $rng1.Start
$ePos = $start1 + 12  # length of "Arbeitspaket"
$delRng = $d.Range($ePos, $ePos + 1)
if ($delRng.Text -ne "e") {
    throw "Unexpected char at edit 1 split point: '$($delRng.Text)'"
}
$delRng.Text = ""
$bmRng1 = $d.Range($ePos, $ePos)
$d.Bookmarks.Add("_GoBack", $bmRng1)
